$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert 3 new rows before row 21 (pushes the old row 21 "TOTAL" row, and
# everything below it, down by 3 rows).
$ws.Rows("21:23").Insert()

# Copy the formatting of the row just above the insertion point (row 20)
# into the 3 new rows so they inherit the same look (borders, number
# formats, fonts, etc.) as the rest of the data block.
$ws.Range("A20:J20").Copy()
$ws.Range("A21:J23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 21 ---
$ws.Range("A21").Value = 45224
$ws.Range("B21").Value = "Collective Trade Links Pvt Ltd"
$ws.Range("C21").Value = "I-C-1-23-453322"
$ws.Range("D21").Value = 44500
$ws.Range("E21").Value = 8010
$ws.Range("H21").Formula = "=D21+E21+F21+G21"
$ws.Range("J21").Value = "24AACCC4813C1ZB"

# --- Row 22 ---
$ws.Range("A22").Value = 45226
$ws.Range("B22").Value = "Collective Trade Links Pvt Ltd"
$ws.Range("C22").Value = "I-C-1-23-453414"
$ws.Range("D22").Value = 22000
$ws.Range("E22").Value = 3960
$ws.Range("H22").Formula = "=D22+E22+F22+G22"
$ws.Range("J22").Value = "24AACCC4813C1ZB"

# --- Row 23 ---
$ws.Range("A23").Value = 45230
$ws.Range("B23").Value = "Namrata Rubber Product Private Limited"
$ws.Range("C23").Value = "64/23-24"
$ws.Range("D23").Value = 70125
$ws.Range("E23").Value = 12622.5
$ws.Range("H23").Value = 82748
$ws.Range("J23").Value = "27AAICN6069P1ZL"

# --- Update the TOTAL row (now row 24) so its SUM ranges include the 3
# newly-added rows. ---
$ws.Range("D24").Formula = "=SUM(D9:D23)"
$ws.Range("E24").Formula = "=SUM(E9:E23)"
$ws.Range("F24").Formula = "=SUM(F9:F23)"
$ws.Range("G24").Formula = "=SUM(G9:G23)"
$ws.Range("H24").Formula = "=SUM(H9:H23)"

# Update the selection to match the saved workbook state.
$ws.Range("C11").Select()
